$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.959.11"
$ws.Range("E2").Value = "  +1.15%  "
$ws.Range("D3").Value = "2.637.73"
$ws.Range("E3").Value = "  +1.86%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "'529.39"
$ws.Range("E5").Value = "  +4.06%  "
$ws.Range("D6").Value = "'155.40"
$ws.Range("E6").Value = "  +1.26%  "
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("E8").Value = "  +0.28%  "
$ws.Range("D9").Value = "'6.67"
$ws.Range("E9").Value = "  -0.76%  "
$ws.Range("E10").Value = "  +5.76%  "
$ws.Range("E11").Value = "  +1.63%  "
$ws.Range("E12").Value = "  +0.15%  "
$ws.Range("D13").Value = "3.095.07"
$ws.Range("E13").Value = "  +1.64%  "
$ws.Range("D14").Value = "60.982.38"
$ws.Range("E14").Value = "  +1.22%  "
$ws.Range("D15").Value = "'21.98"
$ws.Range("E15").Value = "  +2.18%  "
$ws.Range("E16").Value = "  +3.51%  "
$ws.Range("D17").Value = "2.637.35"
$ws.Range("E17").Value = "  +1.85%  "
$ws.Range("E18").Value = "  +0.63%  "
$ws.Range("D19").Value = "'354.10"
$ws.Range("E19").Value = "  +0.47%  "
$ws.Range("D20").Value = "'10.64"
$ws.Range("E20").Value = "  +1.28%  "
$ws.Range("E21").Value = "  +2.24%  "
$ws.Range("E22").Value = "  +0.26%  "
$ws.Range("D23").Value = "'61.71"
$ws.Range("E23").Value = "  +2.32%  "
$ws.Range("E24").Value = "  +2.55%  "
$ws.Range("E25").Value = "  +1.74%  "
$ws.Range("E26").Value = "  +0.17%  "
$ws.Range("D27").Value = "0.0₃0868"
$ws.Range("E27").Value = "  +3.79%  "
$ws.Range("D28").Value = "'7.43"
$ws.Range("E28").Value = "  +1.42%  "
$ws.Range("D29").Value = "'0.999"
$ws.Range("E29").Value = "  -0.10%  "
$ws.Range("D30").Value = "'6.15"
$ws.Range("E30").Value = "  +7.26%  "
$ws.Range("D31").Value = "'19.52"
$ws.Range("E31").Value = "  +0.83%  "
$ws.Range("E32").Value = "  +4.20%  "
$ws.Range("D33").Value = "'150.86"
$ws.Range("E33").Value = "  -0.47%  "
$ws.Range("D34").Value = "'4.17"
$ws.Range("E34").Value = "  +4.68%  "
$ws.Range("B35").Value = "Fetch.AI"
$ws.Range("C35").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D35").Value = "'0.936"
$ws.Range("E35").Value = "  +11.76%  "
$ws.Range("B36").Value = "ImmutableX"
$ws.Range("C36").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D36").Value = "'1.20"
$ws.Range("E36").Value = "  +2.12%  "
$ws.Range("D37").Value = "'0.887"
$ws.Range("E38").Value = "  +1.57%  "
$ws.Range("E39").Value = "  +2.28%  "
$ws.Range("D40").Value = "'306.18"
$ws.Range("E40").Value = "  +3.62%  "
$ws.Range("D41").Value = "'0.643"
$ws.Range("E41").Value = "  +4.29%  "
$ws.Range("E42").Value = "  +2.06%  "
$ws.Range("D43").Value = "'0.0563"
$ws.Range("E43").Value = "  +2.35%  "
$ws.Range("D44").Value = "'0.998"
$ws.Range("E44").Value = "  +0.06%  "
$ws.Range("E45").Value = "  +5.31%  "
$ws.Range("D46").Value = "'19.77"
$ws.Range("E46").Value = "  +0.62%  "
$ws.Range("E47").Value = "  +2.69%  "
$ws.Range("D48").Value = "'19.31"
$ws.Range("E48").Value = "  +8.53%  "
$ws.Range("D49").Value = "'10.33"
$ws.Range("E49").Value = "  +0.25%  "
$ws.Range("D50").Value = "1.988.76"
$ws.Range("E50").Value = "  +0.01%  "
$ws.Range("E51").Value = "  +3.19%  "
